# Update the "想去人数" (F column) figures across the relevant sheets.
# These values were bumped slightly (e.g. refreshed scrape counts) while
# everything else on the rows stays the same.

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1335
$ws1.Range("F4").Value = 1139
$ws1.Range("F5").Value = 1028
$ws1.Range("F6").Value = 1808
$ws1.Range("F7").Value = 573
$ws1.Range("F8").Value = 1210
$ws1.Range("F12").Value = 302
$ws1.Range("F13").Value = 76
$ws1.Range("F15").Value = 703
$ws1.Range("F16").Value = 181
$ws1.Range("F17").Value = 106
$ws1.Range("F21").Value = 162
$ws1.Range("F23").Value = 45
$ws1.Range("F25").Value = 163
$ws1.Range("F27").Value = 880
$ws1.Range("F28").Value = 320
$ws1.Range("F29").Value = 163
$ws1.Range("F33").Value = 17

# Sheet: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 259

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1335
$ws4.Range("F5").Value = 1139
$ws4.Range("F6").Value = 1028
$ws4.Range("F7").Value = 1808
$ws4.Range("F8").Value = 573
$ws4.Range("F9").Value = 1210
$ws4.Range("F14").Value = 302
$ws4.Range("F15").Value = 76
$ws4.Range("F17").Value = 703
$ws4.Range("F18").Value = 181
$ws4.Range("F19").Value = 106
$ws4.Range("F27").Value = 259
$ws4.Range("F28").Value = 259
$ws4.Range("F29").Value = 162
$ws4.Range("F31").Value = 45
$ws4.Range("F33").Value = 163
$ws4.Range("F35").Value = 880
$ws4.Range("F36").Value = 320
$ws4.Range("F39").Value = 163
$ws4.Range("F46").Value = 17
